$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("instruct")

# Clear the cells that are being removed from the trend epi outputs table
$ws.Range("L5").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

# Update the active selection to match the new state
$ws.Range("L5").Select()
